$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain text (prices/percentages with dots as separators,
# not real numbers). Force column D to Text format before writing so that
# numeric-looking strings (e.g. "313.72") are not coerced into real numbers,
# then clear the formatting again so cells end up with no explicit style,
# matching the original file layout.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.350.71'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '1.821.97'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '313.72'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4479'
$ws.Range("E7").Value = '  +2.58%  '
$ws.Range("D8").Value = '0.3754'
$ws.Range("D9").Value = '0.07514'
$ws.Range("E9").Value = '  +3.44%  '
$ws.Range("D10").Value = '0.8854'
$ws.Range("E10").Value = '  +4.74%  '
$ws.Range("D11").Value = '21.01'
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = '1.824.82'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").Value = '6.758'
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").Value = '93.90'
$ws.Range("E14").Value = '  +4.97%  '
$ws.Range("D15").Value = '5.402'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '0.07109'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").Value = '0.000008809'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '15.19'
$ws.Range("E20").Value = '  +1.85%  '
$ws.Range("D21").Value = '27.357.62'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '2.057.28'
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("D26").Value = '2.372'
$ws.Range("E26").Value = '  +7.05%  '
$ws.Range("D27").Value = '151.51'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '18.57'
$ws.Range("D29").Value = '5.364'
$ws.Range("E29").Value = '  +2.49%  '
$ws.Range("D30").Value = '117.94'
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("D31").Value = '0.08842'
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("D32").Value = '0.7925'
$ws.Range("E32").Value = '  +7.04%  '
$ws.Range("D33").Value = '1.198'
$ws.Range("D34").Value = '4.510'
$ws.Range("E34").Value = '  +1.73%  '
$ws.Range("D35").Value = '2.922'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '1.113'
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").Value = '0.01996'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("D39").Value = '0.05329'
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("D40").Value = '7.381'
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("D41").Value = '0.5318'
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").Value = '0.1725'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '2.860'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = '2.331'
$ws.Range("E44").Value = '  +20.01%  '
$ws.Range("D45").Value = '8.743'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").Value = '0.5119'
$ws.Range("E46").Value = '  +7.13%  '
$ws.Range("D47").Value = '10.64'
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("D48").Value = '105.87'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = '1.702'
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  +0.58%  '

$dRange.ClearFormats()
